# Add a new "%Trazabilidad" column (D) with weekly traceability data from the
# ISCII reports, plus two new date rows, and two header cell comments
# documenting the data sources.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: D1 = "%Trazabilidad" -------------------------------
# Copy C1's formatting first so D1 picks up the same bold header style,
# then overwrite with the new header text.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "%Trazabilidad"

# --- Weekly traceability values in column D ---------------------------------
# Only a handful of rows carry an actual value (one per week); the rest of
# column D stays blank, matching the source data.
$ws.Range("D7").Value = 74.7
$ws.Range("D14").Value = 75.3
$ws.Range("D21").Value = 75.9

# --- Two new trailing date rows ---------------------------------------------
$ws.Range("A26").Copy($ws.Range("A27"))
$ws.Range("A27").Value = 44124

$ws.Range("A26").Copy($ws.Range("A28"))
$ws.Range("A28").Value = 44125
$ws.Range("D28").Value = 76.2

# --- Source-documentation comments ------------------------------------------
[void]$ws.Range("B1").AddComment("Fuente: informe diario Sanidad")
[void]$ws.Range("D1").AddComment("De informe ISCII")

# Restore the selection to the cell the author left active (D8).
[void]$ws.Range("D8").Select()
